$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "2019-YDL-NFB-TECH_BASSS-1224-1"
    $ws.Cells.Item($r, 5).Value = "Selasa"
    $ws.Cells.Item($r, 6).Value = 43823.41666666666
}

$ws.Cells.Item(12, 1).Value = "2019-YDL-NFB-TECH_SISMA-1223-1"
$ws.Cells.Item(12, 5).Value = "Senin"
$ws.Cells.Item(12, 6).Value = 43822.41666666666
